$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.904.22"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").Value = "1.641.22"
$ws.Range("E3").Value = "  +1.54%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "215.75"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").Value = "0.5083"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "0.2601"
$ws.Range("E8").Value = "  +1.62%  "
$ws.Range("D9").Value = "0.06472"
$ws.Range("E9").Value = "  +2.03%  "
$ws.Range("D10").Value = "20.30"
$ws.Range("E10").Value = "  +5.67%  "
$ws.Range("D11").Value = "0.07805"
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.662.26"
$ws.Range("E12").Value = "  +2.82%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "4.269"
$ws.Range("E13").Value = "  +0.81%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "1.868.19"
$ws.Range("E14").Value = "  +1.60%  "
$ws.Range("E15").Value = "  +2.30%  "
$ws.Range("D16").Value = "0.0₅7718"
$ws.Range("E16").Value = "  +3.11%  "
$ws.Range("D17").Value = "63.59"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").Value = "25.919.47"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("D19").Value = "1.004"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").Value = "194.64"
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("D21").Value = "4.404"
$ws.Range("E21").Value = "  +1.66%  "
$ws.Range("D22").Value = "9.996"
$ws.Range("E22").Value = "  +2.76%  "
$ws.Range("D23").Value = "6.271"
$ws.Range("E23").Value = "  +5.38%  "
$ws.Range("D24").Value = "1.004"
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("D25").Value = "1.763"
$ws.Range("E25").Value = "  -3.91%  "
$ws.Range("D26").Value = "139.18"
$ws.Range("E26").Value = "  -0.88%  "
$ws.Range("D27").Value = "0.1228"
$ws.Range("E27").Value = "  -2.29%  "
$ws.Range("D28").Value = "6.867"
$ws.Range("E28").Value = "  +2.30%  "
$ws.Range("D29").Value = "15.58"
$ws.Range("E29").Value = "  +1.16%  "
$ws.Range("E30").Value = "  +1.21%  "
$ws.Range("D31").Value = "0.04993"
$ws.Range("E31").Value = "  +3.06%  "
$ws.Range("D32").Value = "3.330"
$ws.Range("E32").Value = "  +1.36%  "
$ws.Range("D33").Value = "3.268"
$ws.Range("E33").Value = "  +3.16%  "
$ws.Range("D34").Value = "1.582"
$ws.Range("E34").Value = "  +2.48%  "
$ws.Range("D35").Value = "2.385"
$ws.Range("E35").Value = "  +0.93%  "
$ws.Range("D36").Value = "0.9104"
$ws.Range("E36").Value = "  +2.34%  "
$ws.Range("D37").Value = "2.586"
$ws.Range("E37").Value = "  +2.59%  "
$ws.Range("D38").Value = "0.5545"
$ws.Range("E38").Value = "  +1.38%  "
$ws.Range("D39").Value = "1.129.40"
$ws.Range("E39").Value = "  +0.71%  "
$ws.Range("D40").Value = "0.01577"
$ws.Range("E40").Value = "  +1.38%  "
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("D42").Value = "5.514"
$ws.Range("E42").Value = "  -0.88%  "
$ws.Range("D43").Value = "99.87"
$ws.Range("E43").Value = "  +3.08%  "
$ws.Range("D44").Value = "0.8022"
$ws.Range("E44").Value = "  +1.54%  "
$ws.Range("D45").Value = "0.0₈109"
$ws.Range("E45").Value = "  -2.76%  "
$ws.Range("D46").Value = "55.77"
$ws.Range("E46").Value = "  +2.32%  "
$ws.Range("D47").Value = "0.4241"
$ws.Range("E47").Value = "  -3.66%  "
$ws.Range("D48").Value = "0.05050"
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("E49").Value = "  +2.56%  "
$ws.Range("D50").Value = "1.001"
$ws.Range("E50").Value = "  +0.37%  "
$ws.Range("E51").Value = "  +0.20%  "
